# Score update from 29th july
# Appends a new round (played 2018-07-29) to the "jul18" sheet and makes
# that sheet the active / selected tab (it was previously "feb18").

$wb = $excel.ActiveWorkbook

# --- 1. Make "jul18" the active sheet (was "feb18") -------------------
$ws = $wb.Worksheets.Item("jul18")
$ws.Activate()

# --- 2. Append the new date header row ---------------------------------
$dateRow = 64
$ws.Cells.Item($dateRow, 1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item($dateRow, 1).Value = [DateTime]"2018-07-29"

# --- 3. Append the 18 hole-by-hole score rows ---------------------------
$holes = @(
  @{ Hole = "Hole 1";  Score = 4; Note = "R"; Putts = 1 },
  @{ Hole = "Hole 2";  Score = 4; Note = "";  Putts = 2 },
  @{ Hole = "Hole 3";  Score = 6; Note = "R"; Putts = 2 },
  @{ Hole = "Hole 4";  Score = 5; Note = "L"; Putts = 1 },
  @{ Hole = "Hole 5";  Score = 3; Note = "";  Putts = 2 },
  @{ Hole = "Hole 6";  Score = 6; Note = "R"; Putts = 3 },
  @{ Hole = "Hole 7";  Score = 5; Note = "L"; Putts = 2 },
  @{ Hole = "Hole 8";  Score = 3; Note = "";  Putts = 2 },
  @{ Hole = "Hole 9";  Score = 3; Note = "S"; Putts = 0 },
  @{ Hole = "Hole 10"; Score = 4; Note = "L"; Putts = 2 },
  @{ Hole = "Hole 11"; Score = 4; Note = "";  Putts = 1 },
  @{ Hole = "Hole 12"; Score = 4; Note = "S"; Putts = 2 },
  @{ Hole = "Hole 13"; Score = 5; Note = "S"; Putts = 2 },
  @{ Hole = "Hole 14"; Score = 5; Note = "S"; Putts = 1 },
  @{ Hole = "Hole 15"; Score = 5; Note = "";  Putts = 2 },
  @{ Hole = "Hole 16"; Score = 4; Note = "R"; Putts = 2 },
  @{ Hole = "Hole 17"; Score = 4; Note = "L"; Putts = 2 },
  @{ Hole = "Hole 18"; Score = 4; Note = "S"; Putts = 1 }
)

$firstDataRow = $dateRow + 1
$r = $firstDataRow
foreach ($hole in $holes) {
  $ws.Cells.Item($r, 1).Value = $hole.Hole
  $ws.Cells.Item($r, 2).Value = $hole.Score
  if ($hole.Note -ne "") {
    $ws.Cells.Item($r, 3).Value = $hole.Note
  }
  $ws.Cells.Item($r, 5).Value = $hole.Putts
  $r = $r + 1
}
$lastDataRow = $r - 1

# --- 4. Totals row -------------------------------------------------------
$totalRow = $r
$ws.Cells.Item($totalRow, 2).Formula = "=SUM(B$($firstDataRow):B$($lastDataRow))"
$ws.Cells.Item($totalRow, 5).Formula = "=SUM(E$($firstDataRow):E$($lastDataRow))"

# --- 5. Update the view: select the new totals cell on "jul18", which is
#        now the active/selected tab (previously "feb18" held that). -------
$ws.Activate()
$ws.Range("B$totalRow").Select()
